$wb = $excel.ActiveWorkbook

# --- Sheet "Login": B2 changes from numeric 12345 to text "Padang123" ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("B2").Value = "Padang123"

# --- Sheet "Tag": B2 changes from "#ffffff" to "#170101" ---
$wsTag = $wb.Worksheets.Item("Tag")
$wsTag.Range("B2").Value = "#170101"

# --- Sheet "TagSave": append rows 107-115 ---
$wsTagSave = $wb.Worksheets.Item("TagSave")
$wsTagSave.Range("A107").Value = "Zv7rFuuzGN"
$wsTagSave.Range("A108").Value = "0KbymM3g4f"
$wsTagSave.Range("A109").Value = "utraHQup0f"
$wsTagSave.Range("A110").Value = "zDqhQbgN0j"
$wsTagSave.Range("A111").Value = "DOoi8yNV0b"
$wsTagSave.Range("A112").Value = "zyotYaM1eN"
$wsTagSave.Range("A113").Value = "Z0JgJdJUzwxz1u6mmKjwt9OwUEBTTwhrzVZPZ2JZaBLhUfOZPfBH5lAuqEXrSl2Z1F67fcierFSFloQteSjUhM4NllRlyYthyjh9OcLZSSnzT5rIb33JfXPYrP6BiyAyXdY9A26l0st8ZrUhxWetvWsZn7qlRsWsocEqKKipwQVrabMuddgZPm8FjVPFWWRWAsvGQ7fm1mbps1rBhVNQZxSTeqFwL2ZHpLbgZMPazrJFOT7B9ombQ4pNEJIbjjIxHvnFzHEk89BQXkOXVXKEK8apWDJNDCLVw0hGMXzw8zwM"
$wsTagSave.Range("A114").Value = "Y4eo58PXf54H81wPWbVFLJfssb3gPatGBIk7yAwL55FIrCKWEVMyHHe5lGwcrLJrocENEErLZJo9LnPaM4vtEGwKhl4WPGcbMM1xQ752WPmwEZ7sPnS82VAgBIURNUKBDZaAtEA9pnFBANd3JyIEDeinvdvSYJ872Z6pePcVbAMV4OdwyHVqjmpPpbQVDzEDBaobV2rOyNrKhB50AE6ZdwFW1vr6v9USBJwUc227IR99qOyZYadRythNCl3s1H6fSLDFtBVeQPorDciKHJ8s5rCQChhXYhIyKDVgiWTg1tPN"

# Row 115 holds a purely-numeric-looking value ("9557487046") that must be
# stored as text, not a number. Enter it as a formula that evaluates to the
# text string, then paste the result back over itself as a value-only paste
# so the final cell is a plain text literal (no formula, no style residue).
$wsTagSave.Range("A115").Formula = "=""9557487046"""
$wsTagSave.Range("A115").Copy()
$wsTagSave.Range("A115").PasteSpecial(-4163)

# --- Sheet "CitySave": append row 5 ---
$wsCitySave = $wb.Worksheets.Item("CitySave")
$wsCitySave.Range("A5").Value = "7N"
